$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Calificaciones")
$ws1.Range("H4").Value = 10
$ws1.Range("T4").Value = 9
$ws1.Range("H5").Value = 10
$ws1.Range("T5").Value = 10
$ws1.Range("H6").Value = 10
$ws1.Range("T6").Value = 9
$ws1.Range("H7").Value = 6
$ws1.Range("H8").Value = 9
$ws1.Range("T8").Value = 8
$ws1.Range("H9").Value = 9
$ws1.Range("T9").Value = 8
$ws1.Range("H10").Value = 10
$ws1.Range("T10").Value = 9
$ws1.Range("H11").Value = 6
$ws1.Range("H12").Value = 10
$ws1.Range("T12").Value = 9
$ws1.Range("H13").Value = 10
$ws1.Range("T13").Value = 9
$ws1.Range("H14").Value = 10
$ws1.Range("T14").Value = 9
$ws1.Range("H15").Value = 10
$ws1.Range("T15").Value = 8
$ws1.Range("H16").Value = 10
$ws1.Range("T16").Value = 9
$ws1.Range("H17").Value = 9
$ws1.Range("T17").Value = 8
$ws1.Range("H18").Value = 10
$ws1.Range("T18").Value = 9
$ws1.Range("H19").Value = 10
$ws1.Range("T19").Value = 9
$ws1.Range("H20").Value = 10
$ws1.Range("T20").Value = 8
$ws1.Range("H21").Value = 10
$ws1.Range("T21").Value = 9
$ws1.Range("H22").Value = 10
$ws1.Range("T22").Value = 10
$ws1.Range("H23").Value = 9
$ws1.Range("T23").Value = 8
$ws1.Range("H24").Value = 10
$ws1.Range("T24").Value = 10
$ws1.Range("H25").Value = 10
$ws1.Range("T25").Value = 9
$ws1.Range("H26").Value = 10
$ws1.Range("T26").Value = 9
$ws1.Range("H27").Value = 10
$ws1.Range("T27").Value = 9
$ws1.Range("H28").Value = 6

$ws2 = $wb.Worksheets.Item("Asistencias")
$ws2.Range("H4").Value = 92
$ws2.Range("N4").Value = 92
$ws2.Range("H5").Value = 96
$ws2.Range("N5").Value = 96
$ws2.Range("H6").Value = 88
$ws2.Range("N6").Value = 88
$ws2.Range("H7").Value = 84
$ws2.Range("N7").Value = 84
$ws2.Range("H8").Value = 92
$ws2.Range("N8").Value = 92
$ws2.Range("H9").Value = 96
$ws2.Range("N9").Value = 96
$ws2.Range("H10").Value = 92
$ws2.Range("N10").Value = 92
$ws2.Range("H11").Value = 88
$ws2.Range("N11").Value = 88
$ws2.Range("H12").Value = 92
$ws2.Range("N12").Value = 92
$ws2.Range("H13").Value = 92
$ws2.Range("N13").Value = 92
$ws2.Range("H14").Value = 92
$ws2.Range("N14").Value = 92
$ws2.Range("H15").Value = 88
$ws2.Range("N15").Value = 88
$ws2.Range("H16").Value = 92
$ws2.Range("N16").Value = 92
$ws2.Range("H17").Value = 88
$ws2.Range("N17").Value = 88
$ws2.Range("H18").Value = 92
$ws2.Range("N18").Value = 92
$ws2.Range("H20").Value = 92
$ws2.Range("N20").Value = 92
$ws2.Range("H22").Value = 94
$ws2.Range("N22").Value = 94
$ws2.Range("H23").Value = 84
$ws2.Range("N23").Value = 84
$ws2.Range("H25").Value = 96
$ws2.Range("N25").Value = 96
$ws2.Range("H26").Value = 88
$ws2.Range("N26").Value = 88
$ws2.Range("H27").Value = 86
$ws2.Range("N27").Value = 86

$ws3 = $wb.Worksheets.Item("Totales")
$ws3.Range("H5").Value = 8.5
